$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Vendor" table (A1:F6) gets reordered and a new "Rockwell" vendor row
# is inserted, growing the table to A1:F7:
#
#   before: Siemens, Asus, Schneider, AVM, Synology
#   after : Schneider, AVM, Rockwell(new), Siemens, Asus, Synology
#
# Easiest/most robust way to reach that exact target layout is to wipe the
# existing data rows (keeping row 1 headers untouched) and rewrite every
# row from scratch in the new order, restoring number formats per cell.
# ---------------------------------------------------------------------------

$dateTimeFmt = "yyyy\-mm\-dd\ hh:mm:ss"
$dateFmt     = "yyyy\-mm\-dd"

# Insert one extra row so the data block grows from 5 to 6 rows (A1:F6 -> A1:F7)
$ws.Rows("6:6").Insert()

# Clear old contents of the whole data block, keep number formats as a base
$ws.Range("A2:F7").ClearContents()

# Style mode per row for columns C (Last_update) / D (Next_update):
#   "normal" -> C: datetime format, D: date format   (the common case)
#   "swap"   -> C: date format,     D: datetime format (Synology, kept as-is from source file)
#   "bothdate" -> C: date format,   D: date format     (new Rockwell row)
function Set-VendorRow {
    param($Row, $Vendor, $Intervall, $LastUpdate, $NextUpdate, $VendorClass, $MaxProducts, $StyleMode)

    $ws.Range("A$Row").Value = $Vendor
    $ws.Range("B$Row").Value2 = $Intervall

    if ($StyleMode -eq "swap") {
        $ws.Range("C$Row").NumberFormat = $dateFmt
        $ws.Range("D$Row").NumberFormat = $dateTimeFmt
    } elseif ($StyleMode -eq "bothdate") {
        $ws.Range("C$Row").NumberFormat = $dateFmt
        $ws.Range("D$Row").NumberFormat = $dateFmt
    } else {
        $ws.Range("C$Row").NumberFormat = $dateTimeFmt
        $ws.Range("D$Row").NumberFormat = $dateFmt
    }
    $ws.Range("C$Row").Value2 = $LastUpdate
    $ws.Range("D$Row").Value2 = $NextUpdate

    if ($VendorClass -ne $null) {
        $ws.Range("E$Row").Value = $VendorClass
    }
    if ($MaxProducts -ne $null) {
        $ws.Range("F$Row").Value2 = $MaxProducts
    }
}

Set-VendorRow 2 "Schneider" 0   44901 44901 "SchneiderElectricScraper" 10    "normal"
Set-VendorRow 3 "AVM"       0   44902 44902 "AVMScraper"               $null "normal"
Set-VendorRow 4 "Rockwell"  0   44934 44935 "RockwellScraper"          $null "bothdate"
Set-VendorRow 5 "Siemens"   100 44894 44993 $null                      $null "normal"
Set-VendorRow 6 "Asus"      100 44893 45261 $null                      $null "normal"
Set-VendorRow 7 "Synology"  0   44902 44902 "SynologyScraper"          $null "swap"

# Column widths: C now best-fits like D (same displayed width, both size for
# the "yyyy-mm-dd hh:mm:ss"/date columns), and E gets its own best-fit width
# for the Vendor_class strings (e.g. "SchneiderElectricScraper"). Column D is
# left untouched so it keeps its original best-fit width/flag. C/E are set to
# the same effective character width Excel previously computed for them.
$ws.Columns.Item(3).ColumnWidth = 16.83
$ws.Columns.Item(5).ColumnWidth = 20.5

# Selection moved to C13 in the saved file
$ws.Range("C13").Select() | Out-Null
